$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.642.58"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "3.605.51"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'609.43"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'148.77"
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.489"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'8.07"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "4.214.08"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'30.10"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "3.544.06"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "66.726.74"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "'11.51"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'6.35"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "'429.27"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "'0.622"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").Value = "'79.14"
$ws.Range("D24").Value = "3.746.63"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").Value = "'8.31"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "3.604.36"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "'25.56"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'176.78"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'2.59"
$ws.Range("E44").Value = "  +8.94%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'1.20"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'24.41"
$ws.Range("E47").Value = "  +5.88%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'25.12"
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("D49").Value = "'7.20"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  -1.26%  "
